# Apply the cryptos-list refresh described by the commit:
# prices/volumes updated, and a few rows re-ranked (B/C/D/E swapped).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    # Route through Formula (not Value) so values are stored as literal
    # text and not re-parsed as numbers/dates by Excel.
    $ws.Range($cell).Formula = $value
}

function Set-NumericLookingTextCell($cell, $value) {
    # Strings such as "378.21" or "1.00" parse as valid numbers, so
    # Excel would silently coerce them (dropping the trailing zero,
    # switching to scientific notation, etc). Force text entry by
    # briefly marking the cell as Text, then restore the original
    # (default/"Normal") cell style so no formatting changes remain.
    $ws.Range($cell).NumberFormat = "@"
    $ws.Range($cell).Formula = $value
    $ws.Range($cell).Style = "Normal"
}

Set-TextCell "D2" "51.058.66"
Set-TextCell "E2" "  -0.05%  "

Set-TextCell "D3" "2.942.27"
Set-TextCell "E3" "  +1.39%  "

Set-TextCell "E4" "  +0.09%  "

Set-NumericLookingTextCell "D5" "378.21"
Set-TextCell "E5" "  +2.44%  "

Set-NumericLookingTextCell "D6" "104.07"
Set-TextCell "E6" "  +0.87%  "

Set-NumericLookingTextCell "D7" "0.541"
Set-TextCell "E7" "  -0.08%  "

Set-TextCell "E8" "  +0.12%  "

Set-NumericLookingTextCell "D9" "0.589"
Set-TextCell "E9" "  -0.11%  "

Set-NumericLookingTextCell "D10" "36.65"
Set-TextCell "E10" "  -0.66%  "

Set-TextCell "E11" "  +0.65%  "

Set-NumericLookingTextCell "D12" "0.0837"
Set-TextCell "E12" "  +0.30%  "

Set-TextCell "D13" "3.411.37"
Set-TextCell "E13" "  +1.59%  "

Set-NumericLookingTextCell "D14" "18.25"
Set-TextCell "E14" "  -0.60%  "

Set-NumericLookingTextCell "D15" "7.42"
Set-TextCell "E15" "  +0.76%  "

Set-TextCell "D16" "2.942.55"
Set-TextCell "E16" "  +1.74%  "

Set-NumericLookingTextCell "D17" "0.949"
Set-TextCell "E17" "  +1.03%  "

Set-TextCell "D18" "51.122.91"
Set-TextCell "E18" "  +0.29%  "

Set-NumericLookingTextCell "D19" "3.31"
Set-TextCell "E19" "  +0.89%  "

Set-NumericLookingTextCell "D20" "7.31"
Set-TextCell "E20" "  +0.99%  "

Set-NumericLookingTextCell "D21" "12.80"
Set-TextCell "E21" "  -0.94%  "

Set-TextCell "D22" "0.0₃0954"
Set-TextCell "E22" "  +1.19%  "

Set-NumericLookingTextCell "D23" "68.84"
Set-TextCell "E23" "  +0.81%  "

Set-NumericLookingTextCell "D24" "260.01"
Set-TextCell "E24" "  +0.27%  "

Set-NumericLookingTextCell "D25" "2.79"
Set-TextCell "E25" "  +3.87%  "

Set-NumericLookingTextCell "D26" "0.167"
Set-TextCell "E26" "  -1.63%  "

Set-NumericLookingTextCell "D27" "7.11"
Set-TextCell "E27" "  +17.59%  "

# Row 28: Dai
Set-TextCell "B28" "Dai"
Set-TextCell "C28" "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-NumericLookingTextCell "D28" "1.00"
Set-TextCell "E28" "  -0.01%  "

# Row 29: Filecoin
Set-TextCell "B29" "Filecoin"
Set-TextCell "C29" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-NumericLookingTextCell "D29" "7.38"
Set-TextCell "E29" "  +2.44%  "

# Row 30: EthereumClassic
Set-TextCell "B30" "EthereumClassic"
Set-TextCell "C30" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-NumericLookingTextCell "D30" "25.72"
Set-TextCell "E30" "  -0.27%  "

# Row 31: Hedera
Set-TextCell "B31" "Hedera"
Set-TextCell "C31" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-NumericLookingTextCell "D31" "0.111"
Set-TextCell "E31" "  +8.25%  "

Set-NumericLookingTextCell "D32" "9.76"
Set-TextCell "E32" "  -1.24%  "

Set-NumericLookingTextCell "D33" "34.38"
Set-TextCell "E33" "  -1.97%  "

Set-NumericLookingTextCell "D34" "2.08"
Set-TextCell "E34" "  -2.35%  "

Set-NumericLookingTextCell "D35" "50.89"
Set-TextCell "E35" "  -0.42%  "

Set-NumericLookingTextCell "D36" "0.0443"
Set-TextCell "E36" "  +5.72%  "

Set-TextCell "E37" "  +0.28%  "

Set-NumericLookingTextCell "D38" "3.04"
Set-TextCell "E38" "  -1.10%  "

Set-NumericLookingTextCell "D39" "16.97"
Set-TextCell "E39" "  +0.20%  "

Set-NumericLookingTextCell "D40" "2.56"
Set-TextCell "E40" "  -4.21%  "

# Row 41: Stellar
Set-TextCell "B41" "Stellar"
Set-TextCell "C41" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-NumericLookingTextCell "D41" "0.115"
Set-TextCell "E41" "  +2.04%  "

# Row 42: ARBITRUM
Set-TextCell "B42" "ARBITRUM"
Set-TextCell "C42" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-NumericLookingTextCell "D42" "1.82"
Set-TextCell "E42" "  -1.00%  "

Set-NumericLookingTextCell "D43" "122.31"
Set-TextCell "E43" "  +4.01%  "

Set-NumericLookingTextCell "D44" "21.89"
Set-TextCell "E44" "  -0.86%  "

Set-NumericLookingTextCell "D45" "0.279"
Set-TextCell "E45" "  +16.07%  "

Set-TextCell "E46" "  -1.08%  "

Set-TextCell "E47" "  +2.33%  "

Set-TextCell "D48" "2.029.02"
Set-TextCell "E48" "  -0.53%  "

Set-NumericLookingTextCell "D49" "3.18"
Set-TextCell "E49" "  +0.41%  "

Set-NumericLookingTextCell "D50" "0.0345"
Set-TextCell "E50" "  +10.48%  "

# Row 51: TrustWalletToken
Set-TextCell "B51" "TrustWalletToken"
Set-TextCell "C51" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-NumericLookingTextCell "D51" "1.27"
Set-TextCell "E51" "  +0.44%  "
